$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Fred VanVleet"
$ws.Cells.Item(2, 2).Value = "PG"
$ws.Cells.Item(2, 3).Value = "Houston Rockets"

$ws.Cells.Item(3, 1).Value = "Collin Sexton"
$ws.Cells.Item(3, 2).Value = "PG,SG"
$ws.Cells.Item(3, 3).Value = "Utah Jazz"

$ws.Cells.Item(4, 1).Value = "Tyus Jones"
$ws.Cells.Item(4, 2).Value = "PG"
$ws.Cells.Item(4, 3).Value = "Phoenix Suns"

$ws.Cells.Item(5, 1).Value = "Jaden Ivey"
$ws.Cells.Item(5, 2).Value = "PG,SG"
$ws.Cells.Item(5, 3).Value = "Detroit Pistons"

$ws.Cells.Item(6, 1).Value = "Anthony Edwards"
$ws.Cells.Item(6, 2).Value = "SG,SF"
$ws.Cells.Item(6, 3).Value = "Minnesota Timberwolves"

$ws.Cells.Item(7, 1).Value = "Giannis Antetokounmpo"
$ws.Cells.Item(7, 2).Value = "PF,C"
$ws.Cells.Item(7, 3).Value = "Milwaukee Bucks"

$ws.Cells.Item(8, 1).Value = "Jaren Jackson Jr."
$ws.Cells.Item(8, 2).Value = "PF,C"
$ws.Cells.Item(8, 3).Value = "Memphis Grizzlies"

$ws.Cells.Item(9, 1).Value = "Kyle Kuzma"
$ws.Cells.Item(9, 2).Value = "PF"
$ws.Cells.Item(9, 3).Value = "Washington Wizards"

$ws.Cells.Item(10, 1).Value = "Draymond Green"
$ws.Cells.Item(10, 2).Value = "PF,C"
$ws.Cells.Item(10, 3).Value = "Golden State Warriors"

$ws.Cells.Item(11, 1).Value = "Jusuf Nurkic"
$ws.Cells.Item(11, 2).Value = "C"
$ws.Cells.Item(11, 3).Value = "Phoenix Suns"

$ws.Cells.Item(12, 1).Value = "Anfernee Simons"
$ws.Cells.Item(12, 2).Value = "PG,SG"
$ws.Cells.Item(12, 3).Value = "Portland Trail Blazers"

$ws.Cells.Item(13, 1).Value = "James Harden"
$ws.Cells.Item(13, 2).Value = "PG,SG"
$ws.Cells.Item(13, 3).Value = "LA Clippers"

$ws.Cells.Item(14, 1).Value = "Ivica Zubac"
$ws.Cells.Item(14, 2).Value = "C"
$ws.Cells.Item(14, 3).Value = "LA Clippers"

$ws.Cells.Item(15, 1).Value = "Jayson Tatum"
$ws.Cells.Item(15, 2).Value = "SF,PF"
$ws.Cells.Item(15, 3).Value = "Boston Celtics"

$ws.Cells.Item(16, 1).Value = "Andrew Wiggins"
$ws.Cells.Item(16, 2).Value = "SF,PF"
$ws.Cells.Item(16, 3).Value = "Golden State Warriors"
